{"js": "// Replace the two-digit multiplication problems in the table cells\n// with the new values, matching the XML diff exactly.\nconst replacements = [\n  [\"74\u00d740=\", \"30\u00d784=\"],\n  [\"23\u00d724=\", \"16\u00d738=\"],\n  [\"87\u00d718=\", \"15\u00d726=\"],\n  [\"73\u00d775=\", \"94\u00d790=\"],\n  [\"83\u00d778=\", \"33\u00d735=\"],\n  [\"94\u00d717=\", \"43\u00d775=\"],\n  [\"40\u00d779=\", \"43\u00d752=\"],\n  [\"13\u00d791=\", \"55\u00d772=\"],\n  [\"86\u00d760=\", \"34\u00d766=\"],\n  [\"39\u00d728=\", \"31\u00d740=\"],\n  [\"34\u00d786=\", \"81\u00d796=\"],\n  [\"85\u00d795=\", \"63\u00d757=\"],\n  [\"13\u00d788=\", \"23\u00d736=\"],\n  [\"79\u00d750=\", \"51\u00d792=\"],\n  [\"97\u00d730=\", \"35\u00d714=\"],\n  [\"40\u00d767=\", \"45\u00d713=\"],\n  [\"32\u00d795=\", \"52\u00d746=\"],\n  [\"59\u00d766=\", \"31\u00d794=\"],\n  [\"50\u00d711=\", \"18\u00d795=\"],\n  [\"97\u00d794=\", \"86\u00d727=\"],\n  [\"81\u00d769=\", \"15\u00d772=\"],\n  [\"52\u00d743=\", \"77\u00d775=\"],\n  [\"16\u00d776=\", \"92\u00d777=\"],\n  [\"21\u00d758=\", \"89\u00d793=\"],\n  [\"94\u00d759=\", \"70\u00d728=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication expression in the table\n# cells with its new value, per the target diff.\n$pairs = @(\n  @(\"74\u00d740=\", \"30\u00d784=\"),\n  @(\"23\u00d724=\", \"16\u00d738=\"),\n  @(\"87\u00d718=\", \"15\u00d726=\"),\n  @(\"73\u00d775=\", \"94\u00d790=\"),\n  @(\"83\u00d778=\", \"33\u00d735=\"),\n  @(\"94\u00d717=\", \"43\u00d775=\"),\n  @(\"40\u00d779=\", \"43\u00d752=\"),\n  @(\"13\u00d791=\", \"55\u00d772=\"),\n  @(\"86\u00d760=\", \"34\u00d766=\"),\n  @(\"39\u00d728=\", \"31\u00d740=\"),\n  @(\"34\u00d786=\", \"81\u00d796=\"),\n  @(\"85\u00d795=\", \"63\u00d757=\"),\n  @(\"13\u00d788=\", \"23\u00d736=\"),\n  @(\"79\u00d750=\", \"51\u00d792=\"),\n  @(\"97\u00d730=\", \"35\u00d714=\"),\n  @(\"40\u00d767=\", \"45\u00d713=\"),\n  @(\"32\u00d795=\", \"52\u00d746=\"),\n  @(\"59\u00d766=\", \"31\u00d794=\"),\n  @(\"50\u00d711=\", \"18\u00d795=\"),\n  @(\"97\u00d794=\", \"86\u00d727=\"),\n  @(\"81\u00d769=\", \"15\u00d772=\"),\n  @(\"52\u00d743=\", \"77\u00d775=\"),\n  @(\"16\u00d776=\", \"92\u00d777=\"),\n  @(\"21\u00d758=\", \"89\u00d793=\"),\n  @(\"94\u00d759=\", \"70\u00d728=\"),\n)\n\n$d = $word.ActiveDocument\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $range = $d.Content\n  $range.Find.ClearFormatting()\n  $range.Find.Replacement.ClearFormatting()\n  $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
